$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (6) - source data now only spans rows 1-5
$ws.Rows.Item(6).Delete()

# Row 2: freelance -> salary, amount unchanged (30000), date shifts by 1 day
$ws.Range("A2").Value = "salary"
$ws.Range("C2").Value = 45988.22928240741

# Row 3: salary/80000 -> Jhjhk/2222, new date
$ws.Range("A3").Value = "Jhjhk"
$ws.Range("B3").Value = "'2222"
$ws.Range("C3").Value = 45983.22928240741

# Row 4: snacks/500 -> books/1199, new date
$ws.Range("A4").Value = "books"
$ws.Range("B4").Value = "'1199"
$ws.Range("C4").Value = 45888.22928240741

# Row 5: books/1199 -> groceries/2200, new date
$ws.Range("A5").Value = "groceries"
$ws.Range("B5").Value = "'2200"
$ws.Range("C5").Value = 45883.22928240741

# B3:B5 got a quote-prefix style from the leading apostrophe above; restore the
# default (unstyled) look so these cells match B2's plain text formatting.
$ws.Range("B3:B5").Style = "Normal"
